$p = $ppt.ActivePresentation

# Remove the last two slides (slide3.xml and slide4.xml in the original
# package, i.e. presentation slide indexes 3 and 4). Delete from the end
# so indexes of the remaining slides don't shift under us.
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()
